$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 44, shifting existing rows 44-116 down to 45-117.
$ws.Rows(44).Insert()

# Populate the newly inserted row 44 with the new weekly price entry.
$ws.Range("A44").Value = 7
$ws.Range("B44").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C44").Value = "Ñuble"
$ws.Range("D44").Value = 45152
$ws.Range("E44").Value = 16
$ws.Range("F44").Value = 100112013
$ws.Range("G44").Value = "Alcachofa"
$ws.Range("H44").Value = "Argentina(o)"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 80
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("M44").Value = 10000
$ws.Range("N44").Value = "`$/caja 50 unidades"
$ws.Range("O44").Value = "Provincia de Limarí"
$ws.Range("P44").Value = 200
$ws.Range("Q44").Value = 50
$ws.Range("R44").Value = "Hortaliza"
